# Swap columns B and C (the AUT_GBR_841850 and CHN_GBR_841850 series)
# across the whole used range of the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value2
    $cVal = $cCell.Value2

    $bCell.Value2 = $cVal
    $cCell.Value2 = $bVal
}
